# Updates the Price (D) and Volume(1h) (E) columns of the cryptos sheet
# with the latest scraped values. D-column values are forced to text
# (leading apostrophe, i.e. Excel quote-prefix) so numeric-looking
# strings such as "322.60" or "1.740.96" keep their exact original
# text representation instead of being coerced into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.614.36"
$ws.Range("E2").Value = "  -2.44%  "
$ws.Range("D3").Value = "'1.748.11"
$ws.Range("E3").Value = "  -3.33%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "'322.60"
$ws.Range("E5").Value = "  -4.63%  "
$ws.Range("D6").Value = "'0.9985"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "'0.4259"
$ws.Range("E7").Value = "  -8.71%  "
$ws.Range("D8").Value = "'0.3618"
$ws.Range("E8").Value = "  -5.54%  "
$ws.Range("D9").Value = "'45.37"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").Value = "'0.07475"
$ws.Range("E10").Value = "  -1.98%  "
$ws.Range("D11").Value = "'1.117"
$ws.Range("E11").Value = "  -3.75%  "
$ws.Range("D12").Value = "'0.9996"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").Value = "'21.53"
$ws.Range("E13").Value = "  -4.40%  "
$ws.Range("D14").Value = "'6.111"
$ws.Range("E14").Value = "  -3.91%  "
$ws.Range("D15").Value = "'7.222"
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("D16").Value = "'1.740.96"
$ws.Range("E16").Value = "  -3.72%  "
$ws.Range("D17").Value = "'0.00001069"
$ws.Range("E17").Value = "  -2.55%  "
$ws.Range("D18").Value = "'87.68"
$ws.Range("E18").Value = "  +7.17%  "
$ws.Range("D19").Value = "'0.06237"
$ws.Range("E19").Value = "  -7.13%  "
$ws.Range("D20").Value = "'0.9983"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").Value = "'16.95"
$ws.Range("E21").Value = "  -3.63%  "
$ws.Range("D22").Value = "'6.131"
$ws.Range("E22").Value = "  -4.63%  "
$ws.Range("D23").Value = "'0.5250"
$ws.Range("E23").Value = "  -5.07%  "
$ws.Range("D24").Value = "'27.601.65"
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("D25").Value = "'11.65"
$ws.Range("E25").Value = "  -2.24%  "
$ws.Range("D26").Value = "'2.316"
$ws.Range("E26").Value = "  -4.42%  "
$ws.Range("D27").Value = "'20.49"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").Value = "'2.372"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").Value = "'151.66"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").Value = "'1.940.40"
$ws.Range("E30").Value = "  -3.65%  "
$ws.Range("D31").Value = "'1.221"
$ws.Range("E31").Value = "  -3.09%  "
$ws.Range("D32").Value = "'126.75"
$ws.Range("E32").Value = "  -5.08%  "
$ws.Range("D33").Value = "'5.708"
$ws.Range("E33").Value = "  -2.76%  "
$ws.Range("D34").Value = "'0.09152"
$ws.Range("E34").Value = "  -4.26%  "
$ws.Range("D35").Value = "'3.680"
$ws.Range("E35").Value = "  -8.82%  "
$ws.Range("E36").Value = "  +4.91%  "
$ws.Range("D37").Value = "'0.02304"
$ws.Range("E37").Value = "  -2.31%  "
$ws.Range("D38").Value = "'0.2137"
$ws.Range("E38").Value = "  -6.47%  "
$ws.Range("D39").Value = "'5.088"
$ws.Range("E39").Value = "  -3.67%  "
$ws.Range("D40").Value = "'0.06100"
$ws.Range("E40").Value = "  -4.57%  "
$ws.Range("D41").Value = "'0.6443"
$ws.Range("E41").Value = "  -3.09%  "
$ws.Range("D42").Value = "'1.192"
$ws.Range("E42").Value = "  -3.74%  "
$ws.Range("D43").Value = "'1.416"
$ws.Range("E43").Value = "  -5.18%  "
$ws.Range("D44").Value = "'7.937"
$ws.Range("E44").Value = "  -4.85%  "
$ws.Range("D45").Value = "'0.9977"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").Value = "'13.84"
$ws.Range("E46").Value = "  -2.90%  "
$ws.Range("D47").Value = "'3.726"
$ws.Range("E47").Value = "  -3.39%  "
$ws.Range("D48").Value = "'0.5906"
$ws.Range("E48").Value = "  -3.99%  "
$ws.Range("D49").Value = "'125.87"
$ws.Range("E49").Value = "  -3.93%  "
$ws.Range("E50").Value = "  -3.90%  "
$ws.Range("D51").Value = "'0.06876"
$ws.Range("E51").Value = "  -3.83%  "
